$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a plain text value (used when the text cannot be mistaken for a number).
function Set-PlainText($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# Helper: set a text value that LOOKS like a number (e.g. "1.00", "589.72") without
# Excel auto-converting it to a numeric cell and without leaving the cells style
# (s attribute) changed from its original (unstyled) state. We briefly mark the cell
# as Text (@) so the assignment is stored verbatim, then restore the original style by
# copying it from a neighboring untouched data cell in the same column (D7, default style).
function Set-NumericLookingText($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $donor = $ws.Range("D7")
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $donor.Style
}

Set-PlainText "D2" "62.507.70"
Set-PlainText "E2" "  +2.54%  "
Set-PlainText "D3" "2.935.83"
Set-PlainText "E3" "  +1.82%  "
Set-NumericLookingText "D4" "1.00"
Set-PlainText "E4" "  -0.03%  "
Set-NumericLookingText "D5" "589.72"
Set-PlainText "E5" "  +0.40%  "
Set-NumericLookingText "D6" "147.23"
Set-PlainText "E6" "  +6.33%  "
Set-PlainText "E7" "  +0.04%  "
Set-PlainText "E8" "  +3.21%  "
Set-PlainText "D9" "2.936.05"
Set-PlainText "E9" "  +1.96%  "
Set-PlainText "E10" "  +4.33%  "
Set-PlainText "E11" "  +9.72%  "
Set-PlainText "E12" "  +2.21%  "
Set-PlainText "E13" "  +8.28%  "
Set-NumericLookingText "D14" "32.27"
Set-PlainText "E14" "  +0.32%  "
Set-PlainText "E15" "  -0.57%  "
Set-PlainText "D16" "3.422.16"
Set-PlainText "E16" "  +1.89%  "
Set-PlainText "D17" "62.473.62"
Set-PlainText "E17" "  +2.64%  "
Set-PlainText "B18" "WrappedEther"
Set-PlainText "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PlainText "D18" "2.939.08"
Set-PlainText "E18" "  +1.83%  "
Set-PlainText "B19" "Polkadot"
Set-PlainText "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-NumericLookingText "D19" "6.63"
Set-PlainText "E19" "  +2.57%  "
Set-NumericLookingText "D20" "433.55"
Set-PlainText "E20" "  +2.41%  "
Set-NumericLookingText "D21" "13.45"
Set-PlainText "E21" "  +1.73%  "
Set-NumericLookingText "D22" "0.663"
Set-PlainText "E22" "  +1.83%  "
Set-PlainText "E23" "  +0.71%  "
Set-NumericLookingText "D24" "80.22"
Set-PlainText "E24" "  +0.59%  "
Set-PlainText "E25" "  +6.18%  "
Set-NumericLookingText "D26" "11.89"
Set-PlainText "E26" "  +5.01%  "
Set-PlainText "E27" "  +2.26%  "
Set-PlainText "E28" "  -0.07%  "
Set-NumericLookingText "D29" "7.28"
Set-PlainText "E29" "  +10.59%  "
Set-PlainText "E30" "  +4.98%  "
Set-PlainText "E31" "  +2.02%  "
Set-PlainText "E32" "  +21.85%  "
Set-PlainText "E33" "  +6.37%  "
Set-NumericLookingText "D34" "26.03"
Set-PlainText "E34" "  +1.90%  "
Set-PlainText "E35" "  -0.01%  "
Set-PlainText "E36" "  +1.82%  "
Set-PlainText "E37" "  +2.97%  "
Set-NumericLookingText "D38" "3.02"
Set-PlainText "E38" "  +8.82%  "
Set-PlainText "E39" "  +1.21%  "
Set-PlainText "E40" "  +5.97%  "
Set-PlainText "E41" "  +1.90%  "
Set-PlainText "E42" "  +0.48%  "
Set-NumericLookingText "D43" "0.274"
Set-PlainText "E43" "  +4.20%  "
Set-NumericLookingText "D44" "39.11"
Set-PlainText "E44" "  +3.93%  "
Set-PlainText "D45" "2.698.33"
Set-PlainText "E45" "  +1.35%  "
Set-NumericLookingText "D46" "135.30"
Set-PlainText "E46" "  +2.73%  "
Set-NumericLookingText "D47" "0.0339"
Set-PlainText "E47" "  +3.50%  "
Set-NumericLookingText "D48" "352.19"
Set-PlainText "E48" "  +1.52%  "
Set-PlainText "E49" "  +0.05%  "
Set-PlainText "E50" "  +2.55%  "
Set-NumericLookingText "D51" "22.49"
